$wb = $excel.ActiveWorkbook

# "Generate Report for Handoff": update the Latest Handoff Datetime for every
# row that is about to be (re-)handed off, i.e. rows whose Status is
# "Handback transform failed" (row 4) or "Ready for handoff" (rows 6-10).
# Rows already "Handed back..." or "In Translation" keep their old datetime.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$zhCnRows = @(4, 6, 7, 8, 9, 10)
foreach ($r in $zhCnRows) {
    $wsZhCn.Range("D$r").Value = "2016-03-01 09:59:47"
}

$wsDeDe = $wb.Worksheets.Item("de-de")
$deDeRows = @(4, 6, 7, 8, 9, 10)
foreach ($r in $deDeRows) {
    $wsDeDe.Range("D$r").Value = "2016-03-01 09:59:58"
}
